$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clear H39 (was "Defense Agency") ---
$ws.Range("H39").Value = ""

# --- 2) Clear H96 (was "Office") ---
$ws.Range("H96").Value = ""

# --- 3) Insert a new row at 135, pushing old rows 135-139 down to 136-140 ---
$ws.Rows.Item(135).Insert()

# Fill in the brand-new row 135 with the "DIRECTOR, DHA" intro line
$ws.Range("A135").Value = "DoDI 6440.02 CH 1.pdf"
$ws.Range("B135").Value = "Clinical Laboratory Improvement Program (CLIP)"

# "2." looks numeric to Excel's auto-detection, so force it through as text
# (pre-format as Text, assign, then drop the number-format again so no
# stray style sticks around on the cell).
$ws.Range("C135").NumberFormat = "@"
$ws.Range("C135").Value = "2."
$ws.Range("C135").ClearFormats()

$ws.Range("D135").Value = "DIRECTOR, DHA. Under the authority, direction, and control of the ASD(HA):"
$ws.Range("E135").Value = "ASDHA;DHA"

$ws.Range("F135").NumberFormat = "@"
$ws.Range("F135").Value = "2."
$ws.Range("F135").ClearFormats()

$ws.Range("G135").Value = "SECRETARIES OF THE MILITARY DEPARTMENTS.  The Secretaries of the Military Departments:"
$ws.Range("H135").Value = "Military Departments"
$ws.Range("I135").Value = "active"

# Update D/E on the shifted rows (136-140) from "SECRETARIES OF..."/"Military Departments"
# to the new "DIRECTOR, DHA"/"ASDHA;DHA" values. F/G/H/I keep their shifted-down content.
for ($r = 136; $r -le 140; $r++) {
    $ws.Range("D" + $r).Value = "DIRECTOR, DHA. Under the authority, direction, and control of the ASD(HA):"
    $ws.Range("E" + $r).Value = "ASDHA;DHA"
}
